$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated results from the power-law model fit.
# Values are written as text so they keep being stored as shared strings,
# matching the workbook's original convention of storing all numbers as text.
$changes = @{
    "B2" = "2.401"
    "D2" = "13"
    "F2" = "200"
    "H2" = "0.037"
    "B3" = "2.425"
    "C3" = "0.132"
    "D3" = "15.107"
    "E3" = "3.96"
    "F3" = "172.879"
    "G3" = "37.793"
    "H3" = "0.038"
    "I3" = "0.01"
    "J3" = "0.49"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    # Drop the temporary text-number-format style so the cell falls back to
    # the workbook's default style, just like every other text cell in the
    # sheet.
    $cell.Style = "Normal"
}
